$wb = $excel.ActiveWorkbook

# --- Locate existing sheets -------------------------------------------------
$wsCypher  = $wb.Worksheets.Item("CypherOutput")
$wsMessage = $wb.Worksheets.Item("Message")

# --- Common message block content (Neo4j connection info) ------------------
$neo4jUrlLabel  = "Neo4j_URL:"
$neo4jUrl       = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userLabel      = "User_name:"
$userName       = "neo4j"
$pwdLabel       = "PWD:"
$pwdValue       = "icdcDBneo4j0"
$cypherLabel    = "Cypher:"
$outputLabel    = "Output:"
$outputPath     = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC34_Canine_Filter_Breed-SaintBnd_Neo4jData.xlsx"

$origCypher = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN [''Saint Bernard''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'

$statCypher = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Saint Bernard'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

function Write-MessageBlock($ws, $startRow, $cypherText) {
    $ws.Cells.Item($startRow + 0, 1).Value = $neo4jUrlLabel
    $ws.Cells.Item($startRow + 1, 1).Value = $neo4jUrl
    $ws.Cells.Item($startRow + 2, 1).Value = $userLabel
    $ws.Cells.Item($startRow + 3, 1).Value = $userName
    $ws.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $ws.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $ws.Cells.Item($startRow + 6, 1).Value = $cypherLabel
    $ws.Cells.Item($startRow + 7, 1).Value = $cypherText
    $ws.Cells.Item($startRow + 8, 1).Value = $outputLabel
    $ws.Cells.Item($startRow + 9, 1).Value = $outputPath
}

# --- New sheet: CypherOutput_Message (placed right after Message) ----------
$wsCypherMsg = $wb.Worksheets.Add($null, $wsMessage)
$wsCypherMsg.Name = "CypherOutput_Message"
Write-MessageBlock $wsCypherMsg 1 $origCypher

# --- New sheet: StatOutput (placed right after CypherOutput_Message) -------
$wsStat = $wb.Worksheets.Add($null, $wsCypherMsg)
$wsStat.Name = "StatOutput"

$wsStat.Cells.Item(1, 1).Value = "number_of_files"
$wsStat.Cells.Item(1, 2).Value = "number_of_sample"
$wsStat.Cells.Item(1, 3).Value = "number_of_cases"
$wsStat.Cells.Item(1, 4).Value = "number_of_study"

# Numeric-looking values must land in the sheet as *text* shared strings
# (matching the source data, which stores "2"/"4"/"3" as strings, not
# numbers). Writing a text formula and then pasting-as-values converts the
# formula's string result into a literal text cell without leaving behind
# any extra cell style.
$wsStat.Cells.Item(2, 1).Formula = "=""2"""
$wsStat.Cells.Item(2, 2).Formula = "=""4"""
$wsStat.Cells.Item(2, 3).Formula = "=""3"""
$wsStat.Cells.Item(2, 4).Formula = "=""2"""
$wsStat.Range("A2:D2").Copy()
$wsStat.Range("A2:D2").PasteSpecial(-4163)

# --- New sheet: StatOutput_Message (placed right after StatOutput) ---------
$wsStatMsg = $wb.Worksheets.Add($null, $wsStat)
$wsStatMsg.Name = "StatOutput_Message"
Write-MessageBlock $wsStatMsg 1 $origCypher
Write-MessageBlock $wsStatMsg 11 $statCypher
